$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("2023-12-08 19:12:36", 0.001),
    @("2023-12-08 19:13:38", 0.004200000000000001),
    @("2023-12-08 19:14:28", 0.003800000000000001),
    @("2023-12-08 19:14:34", 0.0004),
    @("2023-12-08 19:14:43", 0.0002)
)

$startRow = 125
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = [string]$data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
